# Week 13 logging update
# Updates the "R" row (row 3) target depth data on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 373
$wsOff.Range("C3").Value = 283
$wsOff.Range("D3").Value = 88
$wsOff.Range("E3").Value = 45

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 379
$wsDef.Range("C3").Value = 264
$wsDef.Range("D3").Value = 94
$wsDef.Range("E3").Value = 44
$wsDef.Range("F3").Value = 8
